$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Sponsored
Launching 2 BHK in Thane West | Homes in Pokhran Rd by...
raymondtenxera.com
https://www.raymondtenxera.com › official-site › brand
New Launch Homes by Raymond with 38 Habitable Floors, 26500 SqFt Clubhouse, 40+ Amenities. Experience a futuristic lifestyle with Raymond Realty''s Spacious 2 BHK homes...'
$ws.Cells.Item(2, 2).Value = 'Raymond Limited'
$ws.Cells.Item(2, 3).Value = 'India'

$ws.Cells.Item(3, 1).Value = 'Sponsored
Visit Address By GS Thane | Consult an expert & visit site
raymonds-addressbygs.com
https://www.raymonds-addressbygs.com
Bookings Open Addres By GS Thane 6.1 Acre, 2/3/4 Bhk 1.30 Cr Ask Expert & Visit Site. Booking Open For Limited...
Price ₹ / BHK/ Area · View Pricing · Site & Layout Plan · Location Advantage'
$ws.Cells.Item(3, 2).Value = 'Home Bazaar Services Pvt Ltd'
$ws.Cells.Item(3, 3).Value = 'India'

$ws.Cells.Item(4, 1).Value = 'Sponsored
Download Brochure | Launching Final Tower Asteria
luxeoffplans.com
https://www.luxeoffplans.com
Spacious 2, 3 & 4 BHK Homes at Thane. Sample Flat Ready. Download Brochure. Explore...'
$ws.Cells.Item(4, 2).Value = 'RABS NET SOLUTIONS PVT LTD'
$ws.Cells.Item(4, 3).Value = 'India'

$ws.Cells.Item(5, 1).Value = 'Sponsored
OC Received RTMI Homes @0 GST* - Grand Clubhouse...
tenxhabitatraymondrealty.com
https://www.tenxhabitatraymondrealty.com › official_site › book_now
4,200 Sq.Ft. Multipurpose Hall | 28 Seater Mini Theatre | 2,400 Sq.Ft. Fully Equipped Gym'
$ws.Cells.Item(5, 2).Value = 'Raymond Limited'
$ws.Cells.Item(5, 3).Value = 'India'

$ws.Cells.Item(6, 1).Value = 'Sponsored
TenX Habitat Thane Launch | 2, 3 & 4 BHK Starts @ 1.41Cr*
homesfy-property.co.in
https://www.homesfy-property.co.in › tenx › thane
326,500 sq ft clubhouse | 2 and 3 BHK residences near Viviana Mall @ 1.41 CR* Launching Tenx Habitat at Thane. Lavish 2, 3 & 4 BHK Starting At 1.41Cr*.
Price Plan · Our Price · Browse Prices · Floor Plans · View Gallery · Our Gallery'
$ws.Cells.Item(6, 2).Value = 'Homesfy Realty Limited'
$ws.Cells.Item(6, 3).Value = 'India'

$ws.Cells.Item(7, 1).Value = 'Sponsored
RTMI Flats with Zero GST* - Grand Clubhouse @Raymond TenX
tenxhabitatraymondrealty.com
https://www.tenxhabitatraymondrealty.com › official_site › book_now
4,200 Sq.Ft. Multipurpose Hall | 28 Seater Mini Theatre | 2,400 Sq.Ft. Fully Equipped Gym
Real Estate Builders & Construction Company · Thane · Open ⋅ Closes 6 pm
Call us'
$ws.Cells.Item(7, 2).Value = 'Raymond Limited'
$ws.Cells.Item(7, 3).Value = 'India'

$ws.Cells.Item(8, 1).Value = 'Sponsored
Visit Address Tower Thane - Consult an expert & visit site
raymonds-addressbygs.com
https://www.raymonds-addressbygs.com
Bookings Open Addres Tower Thane 6.1 Acre, 2/3/4 Bhk 1.30 Cr Ask Expert & Visit Site'
$ws.Cells.Item(8, 2).Value = 'Home Bazaar Services Pvt Ltd'
$ws.Cells.Item(8, 3).Value = 'India'

$ws.Cells.Item(9, 1).Value = 'Sponsored
Launching Final Tower Asteria | 2, 3 & 4 BHK Apartments
luxeoffplans.com
https://www.luxeoffplans.com
Spacious 2, 3 & 4 BHK Homes at Thane. Sample Flat Ready. Download Brochure. Explore...'
$ws.Cells.Item(9, 2).Value = 'RABS NET SOLUTIONS PVT LTD'
$ws.Cells.Item(9, 3).Value = 'India'

$ws.Cells.Item(10, 1).Value = 'Sponsored
Dynamix Group - Avanya - Dynamix Group Dahisar Project
dynamixavanya.co.in
https://www.dynamixavanya.co.in › dahisar › avanya
Explore Avanya by Dynamix, top builders in Mumbai. Luxury living awaits. Learn more today.'
$ws.Cells.Item(10, 2).Value = 'ANILINE PROPERTIES PRIVATE LIMITED'
$ws.Cells.Item(10, 3).Value = 'India'

$ws.Cells.Item(11, 1).Value = 'Sponsored
2 Bed Houses by Raymond Realty | Pay 20% & Nothing till Jan''25
raymondtenxera.com
https://www.raymondtenxera.com › official-site › brand
Experience a futuristic lifestyle with Raymond Realty''s Spacious 2 BHK homes in Thane West. New Launch Homes by Raymond with...
Location Map · Location Advantages · Contact Us · Get In Touch · Configuration · Overview
Call us'
$ws.Cells.Item(11, 2).Value = 'Raymond Limited'
$ws.Cells.Item(11, 3).Value = 'India'
